# Update EPEX spot prices workbook:
#  - "Prix Spot" sheet: insert a new date column ("14-dec") before the
#    "01-oct." column (column EQ), shifting all subsequent date columns
#    one column to the right, and fill the new column's data rows with "-".
#  - "Gaz" and "CO2" sheets: append a new trailing data row for 2025-12-12.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "Prix Spot"
# ---------------------------------------------------------------------
$wsPrix = $wb.Worksheets.Item("Prix Spot")

# Insert a new column at EQ (column 147), pushing the "01-oct." column
# (and everything after it) one column to the right.
$wsPrix.Range("EQ1").EntireColumn.Insert()

# New header cell for the inserted column.
$wsPrix.Range("EQ1").Value = "14-dec"

# Fill the newly inserted column's data cells (rows 2-25) with "-".
$lastRow = 25
for ($r = 2; $r -le $lastRow; $r++) {
    $wsPrix.Cells.Item($r, 147).Value = "-"
}

# ---------------------------------------------------------------------
# Sheet 2: "Gaz" - append row 177
# ---------------------------------------------------------------------
$wsGaz = $wb.Worksheets.Item("Gaz")
# Force text so the date-like string isn't auto-converted to a date serial,
# then drop back to the default style so the cell carries no formatting
# (matching the rest of the date column).
$wsGaz.Range("A177").NumberFormat = "@"
$wsGaz.Range("A177").Value = "2025-12-12"
$wsGaz.Range("A177").Style = "Normal"
$wsGaz.Range("B177").Value = 26.195

# ---------------------------------------------------------------------
# Sheet 3: "CO2" - append row 177
# ---------------------------------------------------------------------
$wsCo2 = $wb.Worksheets.Item("CO2")
$wsCo2.Range("A177").NumberFormat = "@"
$wsCo2.Range("A177").Value = "2025-12-12"
$wsCo2.Range("A177").Style = "Normal"
$wsCo2.Range("B177").Value = 84.09999999999999
